# MLFlow run triggered automatically
# This edit performs a cyclic re-ordering of several rows in the training
# data sheet (rows 110-113/115-116 form one cycle, rows 118-122 form
# another). Each target row ends up containing the full set of values
# (all 47 columns, A:AU) that used to belong to a different row, while
# rows 109, 114 and 117 stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row -> row whose original values should be copied into it.
# (derived from the canonical OOXML diff)
$mapping = @{
    110 = 116
    111 = 115
    112 = 110
    113 = 111
    115 = 112
    116 = 113
    118 = 121
    119 = 120
    120 = 118
    121 = 122
    122 = 119
}

# Snapshot every source row's full row of values (columns A:AU) before any
# writes happen, so the permutation doesn't clobber data it still needs.
$snapshots = @{}
foreach ($srcRow in $mapping.Values | Sort-Object -Unique) {
    $rng = $ws.Range("A$($srcRow):AU$($srcRow)")
    $snapshots[$srcRow] = $rng.Value2
}

# Now write the snapshot for each target row.
foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $destRng = $ws.Range("A$($targetRow):AU$($targetRow)")
    $destRng.Value2 = $snapshots[$srcRow]
}
